$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update Agosto (row 8) figures
$ws.Range("B8").Value = 817
$ws.Range("C8").Value = 4232

# Update Setembro (row 9) figures
$ws.Range("B9").Value = 1699
$ws.Range("C9").Value = 3741

# Move the active cell selection from J10 to I10
$ws.Activate()
$ws.Range("I10").Select()
